$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Capture all of the pre-existing data first (labels + word counts) so that
# we are free to lay out the new, wider table without worrying about
# clobbering a cell before we have read it.
# ---------------------------------------------------------------------------
$labelAbstract     = $ws.Range("A1").Value2
$labelIntroduction  = $ws.Range("A2").Value2
$labelBox           = $ws.Range("A3").Value2
$labelDeriving      = $ws.Range("A4").Value2
$labelWhen          = $ws.Range("A5").Value2
$labelReindeer      = $ws.Range("A6").Value2
$labelTools         = $ws.Range("A7").Value2
$labelDiscussion    = $ws.Range("A8").Value2
$labelOther         = $ws.Range("A9").Value2
$labelRefs          = $ws.Range("A10").Value2
$labelCaptions      = $ws.Range("A11").Value2
$labelTotal         = $ws.Range("A13").Value2
$labelMaxAllowed    = $ws.Range("B18").Value2

$vAbstract     = $ws.Range("B1").Value2
$vIntroduction = $ws.Range("B2").Value2
$vBox          = $ws.Range("B3").Value2
$vDeriving     = $ws.Range("B4").Value2
$vWhen         = $ws.Range("B5").Value2
$vReindeer     = $ws.Range("B6").Value2
$vTools        = $ws.Range("B7").Value2
$vDiscussion   = $ws.Range("B8").Value2
$vOther        = $ws.Range("B9").Value2
$vRefs         = $ws.Range("B10").Value2
$vCaptions     = $ws.Range("B11").Value2

# ---------------------------------------------------------------------------
# Wipe the sheet clean; we'll rebuild every cell explicitly below. This also
# resets the shared-string table so we can control the order new strings are
# registered in.
# ---------------------------------------------------------------------------
$ws.Cells.ClearContents()

# ---------------------------------------------------------------------------
# Re-insert the labels that already existed, in their original relative
# order, which reproduces shared-string indices 0-11 exactly as before.
# ---------------------------------------------------------------------------
$ws.Range("A3").Value  = $labelAbstract
$ws.Range("A4").Value  = $labelIntroduction
$ws.Range("A5").Value  = $labelBox
$ws.Range("A6").Value  = $labelDeriving
$ws.Range("A8").Value  = $labelWhen
$ws.Range("A9").Value  = $labelReindeer
$ws.Range("A10").Value = $labelTools
$ws.Range("A11").Value = $labelDiscussion
$ws.Range("A12").Value = $labelOther
$ws.Range("A13").Value = $labelRefs
$ws.Range("A14").Value = $labelCaptions
$ws.Range("A16").Value = $labelTotal
$ws.Range("B20").Value = $labelMaxAllowed

# ---------------------------------------------------------------------------
# Register the brand-new strings, in the order the original author first
# typed them (reconstructed from the saved shared-string table), so freshly
# minted shared-string indices land exactly where the target file expects.
# ---------------------------------------------------------------------------
$ws.Range("B1").Value  = "V1"
$ws.Range("C1").Value  = "V2_after_bram_review"
$ws.Range("B21").Value = "Max abstract: 350"
$ws.Range("A2").Value  = "Title page"
$ws.Range("D1").Value  = "v3_before_Manu"
$ws.Range("A7").Value  = "Estimating"
$ws.Range("E1").Value  = "v4_after_Manu"
$ws.Range("F1").Value  = "v5_after_Audun"

# ---------------------------------------------------------------------------
# Fill in the numeric word-count grid (columns B..F) for every section.
# ---------------------------------------------------------------------------
# Title page (row 2) - only counted from V2 onward
$ws.Range("C2").Value = 98
$ws.Range("D2").Value = 111
$ws.Range("E2").Value = 105
$ws.Range("F2").Value = 105

# Abstract (row 3)
$ws.Range("B3").Value = $vAbstract
$ws.Range("C3").Value = 358
$ws.Range("D3").Value = 352
$ws.Range("E3").Value = 372
$ws.Range("F3").Value = 386

# Introduction (row 4)
$ws.Range("B4").Value = $vIntroduction
$ws.Range("C4").Value = 1039
$ws.Range("D4").Value = 1125
$ws.Range("E4").Value = 1110
$ws.Range("F4").Value = 1100

# Box (row 5)
$ws.Range("B5").Value = $vBox
$ws.Range("C5").Value = 434
$ws.Range("D5").Value = 531
$ws.Range("E5").Value = 502
$ws.Range("F5").Value = 470

# Deriving (row 6)
$ws.Range("B6").Value = $vDeriving
$ws.Range("C6").Value = 825
$ws.Range("D6").Value = 958
$ws.Range("E6").Value = 1023
$ws.Range("F6").Value = 904

# Estimating (row 7) - brand-new section
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 207
$ws.Range("E7").Value = 226
$ws.Range("F7").Value = 223

# When (row 8)
$ws.Range("B8").Value = $vWhen
$ws.Range("C8").Value = 535
$ws.Range("D8").Value = 463
$ws.Range("E8").Value = 429
$ws.Range("F8").Value = 357

# Reindeer (row 9)
$ws.Range("B9").Value = $vReindeer
$ws.Range("C9").Value = 1040
$ws.Range("D9").Value = 995
$ws.Range("E9").Value = 952
$ws.Range("F9").Value = 954

# Tools (row 10)
$ws.Range("B10").Value = $vTools
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0

# Discussion (row 11)
$ws.Range("B11").Value = $vDiscussion
$ws.Range("C11").Value = 1714
$ws.Range("D11").Value = 1643
$ws.Range("E11").Value = 1570
$ws.Range("F11").Value = 1493

# Other (row 12) - no F value
$ws.Range("B12").Value = $vOther
$ws.Range("C12").Value = 140
$ws.Range("D12").Value = 282
$ws.Range("E12").Value = 299

# Refs (row 13)
$ws.Range("B13").Value = $vRefs
$ws.Range("C13").Value = 1239
$ws.Range("D13").Value = 1250
$ws.Range("E13").Value = 1297
$ws.Range("F13").Value = 1297

# captions (row 14)
$ws.Range("B14").Value = $vCaptions
$ws.Range("C14").Value = 428
$ws.Range("D14").Value = 473
$ws.Range("E14").Formula = "=134+76+119+80+51"
$ws.Range("F14").Formula = "=134+76+119+80+51"

# total: (row 16)
$ws.Range("B16:C16").Formula = "=SUM(B2:B11,B13:B14)"
$ws.Range("D16").Formula = "=SUM(D2:D11,D13:D14)"
$ws.Range("E16").Formula = "=SUM(E2:E11,E13:E14)"
$ws.Range("F16").Formula = "=SUM(F2:F11,F13:F14)"

# ---------------------------------------------------------------------------
# Column widths / view tweaks
# ---------------------------------------------------------------------------
$ws.Columns("C").ColumnWidth = 19.1666666666667
$ws.Columns("D").ColumnWidth = 14.3333333333333

$ws.Range("E16:F16").Select()

Write-Host "done"
